$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 498
$col = 3  # Column C ("Förändrad")

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    if ($cell.Value2 -eq 45179) {
        $cell.Value2 = 45180
    }
}
